# Generate Report for Handoff
# Updates the "b.md" row (row 3) across the Overview, zh-cn and de-de sheets
# to reflect that the file is now "Ready for handoff" instead of
# "Handed back: in sync with en-US", with a freshly generated handoff
# xliff file, timestamp, and an error detail describing a stale handback.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-02 00:47:05"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces Excel to store this as text instead of a boolean;
# re-apply the plain style afterwards so no stray quote-prefix formatting sticks.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = $wsZhCn.Range("E3").Style
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-02 00:46:58"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/057cc89f20635344c1f4e05ae978f23430b8b272/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/853ad9f54b05272715e0671bf7243953f9ebfab0/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces Excel to store this as text instead of a boolean;
# re-apply the plain style afterwards so no stray quote-prefix formatting sticks.
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = $wsDeDe.Range("E3").Style
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-02 00:47:05"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/057cc89f20635344c1f4e05ae978f23430b8b272/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/853ad9f54b05272715e0671bf7243953f9ebfab0/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
